$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.848522901535034
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.275805950164795
$ws.Range("D1").Value = 1.475088357925415
$ws.Range("E1").Value = 1.204766154289246
